$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search_Premium")

$src = $ws.Range("A2:I4")
$src.Copy()
$ws.Range("A5").PasteSpecial(-4104)
$src.Copy()
$ws.Range("A8").PasteSpecial(-4104)
